# Add Item Config And So on
# Adds a new "Icon" field row (row 19) to the Property sheet of the Item
# workbook, matching the style/format of the existing rows, and registers
# the accompanying defined name that the authoring environment maintains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 19: Icon field -------------------------------------------------
$ws.Range("A19").Value = "Icon"
$ws.Range("A19").NumberFormat = "@"

$ws.Range("B19").Value = "string"
$ws.Range("B19").NumberFormat = "@"

$ws.Range("C19").Value = $false
$ws.Range("D19").Value = $false
$ws.Range("E19").Value = $false
$ws.Range("F19").Value = $true

$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0

$ws.Range("I19").Value = "Friend"
$ws.Range("I19").NumberFormat = "@"

$ws.Range("J19").Value = "物品显示Icon"
$ws.Range("J19").NumberFormat = "@"

# --- Workbook-level defined name (LOCAL_MYSQL_DATE_FORMAT) -----------------
$localMysqlDateFormat = "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)"
$definedName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $localMysqlDateFormat)
$definedName.Visible = $false

# --- Restore the cursor position left behind by the edit -------------------
$ws.Range("J17").Select() | Out-Null
